# Updated cryptos list values (price + 1h volume change) per the target diff.
# Cells whose new text parses as a plain number (e.g. "0.999", "605.16") are
# written with a leading apostrophe so Excel keeps them as literal text
# (matching the workbook's inlineStr/text-typed cells) instead of silently
# converting them to numeric cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = '70.033.92' }
    @{ Cell = "E2"; Value = '  +0.32%  ' }
    @{ Cell = "D3"; Value = '3.611.06' }
    @{ Cell = "E3"; Value = '  +3.28%  ' }
    @{ Cell = "D4"; Value = '''0.999' }
    @{ Cell = "E4"; Value = '  -0.01%  ' }
    @{ Cell = "D5"; Value = '''605.16' }
    @{ Cell = "E5"; Value = '  +0.48%  ' }
    @{ Cell = "D6"; Value = '''195.58' }
    @{ Cell = "E6"; Value = '  -0.45%  ' }
    @{ Cell = "D7"; Value = '''0.627' }
    @{ Cell = "E9"; Value = '  -1.59%  ' }
    @{ Cell = "D10"; Value = '''0.651' }
    @{ Cell = "E10"; Value = '  -0.13%  ' }
    @{ Cell = "D11"; Value = '''53.98' }
    @{ Cell = "E11"; Value = '  -0.12%  ' }
    @{ Cell = "D12"; Value = '''0.0000305' }
    @{ Cell = "D13"; Value = '''9.56' }
    @{ Cell = "E13"; Value = '  -0.02%  ' }
    @{ Cell = "D14"; Value = '4.190.61' }
    @{ Cell = "E14"; Value = '  +3.33%  ' }
    @{ Cell = "D15"; Value = '''13.21' }
    @{ Cell = "E15"; Value = '  +4.97%  ' }
    @{ Cell = "D16"; Value = '''591.60' }
    @{ Cell = "E16"; Value = '  -1.60%  ' }
    @{ Cell = "D17"; Value = '''19.21' }
    @{ Cell = "E17"; Value = '  +0.72%  ' }
    @{ Cell = "D18"; Value = '70.251.93' }
    @{ Cell = "E18"; Value = '  +0.45%  ' }
    @{ Cell = "D19"; Value = '3.620.65' }
    @{ Cell = "E19"; Value = '  +3.30%  ' }
    @{ Cell = "E20"; Value = '  +1.57%  ' }
    @{ Cell = "E21"; Value = '  +0.52%  ' }
    @{ Cell = "D22"; Value = '''17.82' }
    @{ Cell = "E22"; Value = '  -2.31%  ' }
    @{ Cell = "D23"; Value = '''5.18' }
    @{ Cell = "E23"; Value = '  +3.08%  ' }
    @{ Cell = "D24"; Value = '''102.71' }
    @{ Cell = "E24"; Value = '  -1.96%  ' }
    @{ Cell = "D25"; Value = '''4.63' }
    @{ Cell = "E25"; Value = '  +1.13%  ' }
    @{ Cell = "D26"; Value = '''3.06' }
    @{ Cell = "E26"; Value = '  -1.26%  ' }
    @{ Cell = "E27"; Value = '  -1.67%  ' }
    @{ Cell = "E28"; Value = '  -1.28%  ' }
    @{ Cell = "D29"; Value = '''33.90' }
    @{ Cell = "E29"; Value = '  +0.95%  ' }
    @{ Cell = "D30"; Value = '''4.42' }
    @{ Cell = "E30"; Value = '  -0.84%  ' }
    @{ Cell = "D31"; Value = '''7.13' }
    @{ Cell = "E31"; Value = '  -0.97%  ' }
    @{ Cell = "D32"; Value = '''12.35' }
    @{ Cell = "E32"; Value = '  -2.72%  ' }
    @{ Cell = "E33"; Value = '  +1.42%  ' }
    @{ Cell = "D35"; Value = '0.0₃0896' }
    @{ Cell = "E35"; Value = '  +11.79%  ' }
    @{ Cell = "D36"; Value = '3.952.72' }
    @{ Cell = "E36"; Value = '  +5.84%  ' }
    @{ Cell = "D37"; Value = '''3.18' }
    @{ Cell = "E37"; Value = '  +6.56%  ' }
    @{ Cell = "D38"; Value = '''529.25' }
    @{ Cell = "E38"; Value = '  +3.80%  ' }
    @{ Cell = "E39"; Value = '  +0.17%  ' }
    @{ Cell = "D40"; Value = '''37.28' }
    @{ Cell = "E40"; Value = '  +1.90%  ' }
    @{ Cell = "E41"; Value = '  +1.05%  ' }
    @{ Cell = "D42"; Value = '''3.55' }
    @{ Cell = "E42"; Value = '  +1.22%  ' }
    @{ Cell = "E43"; Value = '  -1.80%  ' }
    @{ Cell = "D44"; Value = '''0.0456' }
    @{ Cell = "E44"; Value = '  +0.08%  ' }
    @{ Cell = "B45"; Value = 'ApeXProtocol' }
    @{ Cell = "C45"; Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex' }
    @{ Cell = "D45"; Value = '''3.38' }
    @{ Cell = "E45"; Value = '  +2.08%  ' }
    @{ Cell = "B46"; Value = 'ThetaToken' }
    @{ Cell = "C46"; Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta' }
    @{ Cell = "D46"; Value = '''2.87' }
    @{ Cell = "E46"; Value = '  +2.01%  ' }
    @{ Cell = "E47"; Value = '  +0.98%  ' }
    @{ Cell = "D48"; Value = '''8.62' }
    @{ Cell = "E48"; Value = '  -1.30%  ' }
    @{ Cell = "E49"; Value = '  -0.01%  ' }
    @{ Cell = "D50"; Value = '''0.000255' }
    @{ Cell = "E50"; Value = '  +6.39%  ' }
    @{ Cell = "E51"; Value = '  +2.86%  ' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
